$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B37 was stored as text "3" - make it a real number 3
$ws.Range("B37").Value = 3

# Add new row 38 with annotation data
$ws.Range("A38").Value = "Sunsi Wu"
# Force B38 to be stored as text "4" (not a number), matching the source data,
# then drop the number-format override so no stray style sticks to the cell.
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "4"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "well"
$ws.Range("D38").Value = "SMY"
$ws.Range("E38").Value = "MET"
$ws.Range("F38").Value = "d0296b92-10f5-497e-8726-aae675ac805b"
$ws.Range("G38").Value = "rJl3yM-Ab_annotated.xlsx"
$ws.Range("H38").Value = "The new method is motivated well and departs from prior work."
